$wb = $excel.ActiveWorkbook

# Last existing sheet - the two new sheets get added after it
$ws3 = $wb.Worksheets.Item("GradeOneStudentCredentials")

# Add the two new sheets at the end, in order, named for Grade 4 and Grade 9
$ws4 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws3)
$ws4.Name = "Grade4StudentCredentials"

$ws5 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws4)
$ws5.Name = "Grade9StudentCredentials"

# Populate Grade4StudentCredentials (same header layout as the other roster sheets)
$ws4.Range("A1").Value = "UserId"
$ws4.Range("B1").Value = "Password"
$ws4.Range("C1").Value = "UserName"
$ws4.Range("D1").Value = "Signature"

$ws4.Range("A2").Value = "Auto2021_07_16_01_25_46_129"
$ws4.Range("B2").Value = "Password@123"

$ws4.Range("A3").Value = "Auto2021_07_16_01_33_33_468"
$ws4.Range("B3").Value = "Password@123"

$ws4.Rows(1).RowHeight = 16
$ws4.Columns(1).ColumnWidth = 29.666666666666668

# Mirror the header formatting used on the other roster sheets
$ws4.Range("A1").Style = $ws3.Range("A1").Style
$ws4.Range("B1").Style = $ws3.Range("B1").Style

# Populate Grade9StudentCredentials (header row only)
$ws5.Range("A1").Value = "UserId"
$ws5.Range("B1").Value = "Password"
$ws5.Range("C1").Value = "UserName"
$ws5.Range("D1").Value = "Signature"

$ws5.Rows(1).RowHeight = 16

$ws5.Range("A1").Style = $ws3.Range("A1").Style
$ws5.Range("B1").Style = $ws3.Range("B1").Style

# Update selections to match the new authoring state
$ws3.Range("A1:D1").Select() | Out-Null
$ws4.Range("A1:D1").Select() | Out-Null
$ws5.Range("G21").Select() | Out-Null

# Grade9StudentCredentials is the newly active / selected tab
$ws5.Activate() | Out-Null
